$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.271.52"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "1.791.09"
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07057"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8703"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07769"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("E12").Value = "  -3.05%  "
$ws.Range("D13").Value = "1.747.35"
$ws.Range("E13").Value = "  -6.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.260"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.312"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.60%  "
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008479"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "26.332.18"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.64%  "
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.009.41"
$ws.Range("E24").Value = "  -3.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.980"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.025"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.820"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08639"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.024"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.429"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7120"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.88%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.634"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.099"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.006"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05077"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.861"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4906"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("E44").Value = "  -5.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.948"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.82%  "
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.871"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.575"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05939"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.11%  "
